$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data values for rows 2-6, columns A-I (J is the shared string "train_dim2_1" in all rows)
$data = @(
    @(1, 9, 2, 6, 5, -3, 3, 34, 5),
    @(2, 5, 0, 4, 5, -1, 5, 56, 5),
    @(3, 8, 1, 3, 2, -5, 1, 12, 5),
    @(4, 5, 2, 3, 6, -2, 4, 45, 5),
    @(5, 9, 1, 5, 3, -4, 2, 23, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}
